# Cook County Assessment Ratios - add 2017 data, rename sheet/tab

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the workbook title
$ws.Name = "Cook County Assessment Ratios"

# 2017 assessment-district ratios (Cook County Equalization Factor: 2.9627)
$districts2017 = @(
    [PSCustomObject]@{ Name = "Barrington";    Ratio = 10.22 },
    [PSCustomObject]@{ Name = "Elk Grove";     Ratio = 9.03 },
    [PSCustomObject]@{ Name = "Evanston";      Ratio = 8.78 },
    [PSCustomObject]@{ Name = "Hanover";       Ratio = 9.62 },
    [PSCustomObject]@{ Name = "Leyden";        Ratio = 9.47 },
    [PSCustomObject]@{ Name = "Maine";         Ratio = 8.44 },
    [PSCustomObject]@{ Name = "New Trier";     Ratio = 9.07 },
    [PSCustomObject]@{ Name = "Niles";         Ratio = 8.83 },
    [PSCustomObject]@{ Name = "Northfield";    Ratio = 8.97 },
    [PSCustomObject]@{ Name = "Norwood park";  Ratio = 8.53 },
    [PSCustomObject]@{ Name = "Palatine";      Ratio = 9.27 },
    [PSCustomObject]@{ Name = "Schaumburg";    Ratio = 8.93 },
    [PSCustomObject]@{ Name = "Wheeling";      Ratio = 8.64 },
    [PSCustomObject]@{ Name = "Berwyn";        Ratio = 7.05 },
    [PSCustomObject]@{ Name = "Bloom";         Ratio = 9.5 },
    [PSCustomObject]@{ Name = "Bremen";        Ratio = 8.94 },
    [PSCustomObject]@{ Name = "Calumet";       Ratio = 9.25 },
    [PSCustomObject]@{ Name = "Cicero";        Ratio = 6.68 },
    [PSCustomObject]@{ Name = "Lemont";        Ratio = 8.46 },
    [PSCustomObject]@{ Name = "Lyons";         Ratio = 7.96 },
    [PSCustomObject]@{ Name = "Oak Park";      Ratio = 7.74 },
    [PSCustomObject]@{ Name = "Orland";        Ratio = 8.69 },
    [PSCustomObject]@{ Name = "Palos";         Ratio = 8.56 },
    [PSCustomObject]@{ Name = "Proviso";       Ratio = 7.76 },
    [PSCustomObject]@{ Name = "Rich";          Ratio = 9.44 },
    [PSCustomObject]@{ Name = "River Forest";  Ratio = 7.58 },
    [PSCustomObject]@{ Name = "Riverside";     Ratio = 7.69 },
    [PSCustomObject]@{ Name = "Stickney";      Ratio = 8.08 },
    [PSCustomObject]@{ Name = "Thornton";      Ratio = 9.55 },
    [PSCustomObject]@{ Name = "Worth";         Ratio = 8.81 }
)

$equalizationFactor2017 = 2.9627

$row = 332
foreach ($d in $districts2017) {
    $ws.Cells.Item($row, 1).Value = 2017
    $ws.Cells.Item($row, 2).Value = $d.Name
    $ws.Cells.Item($row, 3).Value = $d.Ratio
    $ws.Cells.Item($row, 4).Value = $equalizationFactor2017
    $row++
}

$ws.Range("F355").Select()
